$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

$ws.Range("A14").Value = "OR.0038.0031"
$ws.Range("A3").Copy($ws.Range("A15"))
$ws.Range("A16").Value = "OR.0046.0031"
$ws.Range("A3").Copy($ws.Range("A17"))
